$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 height shrinks from 90 to 75 (wrapped text needs less room now).
$ws.Rows("2").RowHeight = 75

# F2 previously held an (empty) cell with style "1"; it's removed entirely.
$ws.Range("F2").Clear()

# New "Testing 2" test-data values spread across E2, G2 and H2 (plain, unstyled cells),
# plus numeric inspector/capacity figures in I2 and J2.
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").Value = "Testing 2"
$ws.Range("G2").Value = "Testing 2"
$ws.Range("H2").Value = "Testing 2"
$ws.Range("I2").Value = 425
$ws.Range("J2").Value = 525

# Move/save the active selection to F2, matching the saved view state.
[void]$ws.Range("F2").Select()
